$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 must stay a text string ("10,135,145") even though it looks like a
# thousands-grouped number -- force the cell to Text first so Excel doesn't
# silently coerce it into a numeric value.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "10,135,145"

$ws.Range("A2").Value = "H,ABC-123,C,sz,119,13:00:00"
$ws.Range("A3").Value = "A,S-578-Il,B,m,134,18:09:23"
$ws.Range("A4").Value = "GR,IBX-5470,A,b,87,07:10:00"
$ws.Range("A5").Value = "GR,IBX-5470,B,b,40,07:14:00"
$ws.Range("A6").Value = "GR,IBX-5470,C,b,40,07:20:00"
$ws.Range("A7").Value = "H,KZK-252,A,b,101,06:10:12"
$ws.Range("A8").Value = "H,KZK-252,B,b,101,09:50:00"
$ws.Range("A9").Value = "H,KZK-252,C,b,101,10:12:00"

# New rows 5-9 should carry the same explicit row height as the
# pre-existing rows 1-4 (matches sheetFormatPr's default of 29.25).
$ws.Rows.Item(5).RowHeight = 29.25
$ws.Rows.Item(6).RowHeight = 29.25
$ws.Rows.Item(7).RowHeight = 29.25
$ws.Rows.Item(8).RowHeight = 29.25
$ws.Rows.Item(9).RowHeight = 29.25

# Match the saved selection state (active cell B4).
$ws.Range("B4").Select()
